# Update the Reg iExam TC row (row 2) on the active worksheet with the
# newly generated Candidate ID, User Name, Exam Password, First Name and
# Last Name values, and the refreshed numeric Client Id.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 'CLfwZ777'
$ws.Range("B2").Value = 23080256
$ws.Range("C2").Value = 'kaqhund46'
$ws.Range("D2").Value = 'wY8$X%2n'
$ws.Range("F2").Value = 'SIJjvdpw'
$ws.Range("G2").Value = 'YBKM'
